$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds the "K" metric (formerly "Strike#"). Recalculated values.
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
